$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A labels (rewritten per diff) ---
$ws.Range("A1").Value = "Total time taken for the ride"
$ws.Range("A2").Value = "Actual Ampere-hours (Ah)"
$ws.Range("A3").Value = "Actual Watt-hours (Wh)"
$ws.Range("A4").Value = "Starting SoC (Ah)"
$ws.Range("A5").Value = "Ending SoC (Ah)"
$ws.Range("A6").Value = "Starting SoC (%)"
$ws.Range("A7").Value = "Ending SoC (%)"
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("A11").Value = "Mode"
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("A43").Value = "Time spent in 80-90 km/h"

# --- Column B values (rewritten per diff) ---
$ws.Range("B1").Value = 0.03550451388888889
$ws.Range("B2").Value = 32.97032416666666
$ws.Range("B3").Value = 1673.543363610833
$ws.Range("B4").Value = 39.268
$ws.Range("B5").Value = 7.33
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 18
$ws.Range("B8").Value = 30.98331175696487
$ws.Range("B9").Value = 54.01434736022466
$ws.Range("B10").Value = 81
$ws.Range("B11").Value = "Custom mode`n71.07%`nEco mode`n17.72%`nSports mode`n0.07%"
$ws.Range("B12").Value = 5421.410500000001
$ws.Range("B13").Value = -1974.682434939036
$ws.Range("B14").Value = 1.447407298055555
$ws.Range("B15").Value = 0.08641285212993494
$ws.Range("B16").Value = 3.332
$ws.Range("B17").Value = 3.071
$ws.Range("B18").Value = 0.2609999999999997
$ws.Range("B19").Value = 38
$ws.Range("B20").Value = 48
$ws.Range("B21").Value = 10
$ws.Range("B22").Value = 70
$ws.Range("B23").Value = 66
$ws.Range("B24").Value = 65
$ws.Range("B25").Value = 68
$ws.Range("B26").Value = 100
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 48
$ws.Range("B29").Value = 38
$ws.Range("B30").Value = 10
$ws.Range("B31").Value = 53
$ws.Range("B32").Value = 1.747427180833333
$ws.Range("B33").Value = 0.0000001582642449039356
$ws.Range("B34").Value = 14.05714285714286
$ws.Range("B35").Value = 8.607142857142858
$ws.Range("B36").Value = 3.9
$ws.Range("B37").Value = 9.178571428571429
$ws.Range("B38").Value = 18.15
$ws.Range("B39").Value = 12.09642857142857
$ws.Range("B40").Value = 12.00357142857143
$ws.Range("B41").Value = 14.26785714285714
$ws.Range("B42").Value = 7.646428571428572
$ws.Range("B43").Value = 0
